# Auto-generated Excel COM-interop script applying scheduled data refresh
# to the Goblin_Profits-style crafting profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# are refreshed with newly fetched market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 919.3333
$ws.Range("I20").Value = 919.3333
$ws.Range("K20").Value = 919.3333
$ws.Range("M20").Value = -689.3333

$ws.Range("H35").Value = 919.3333
$ws.Range("I35").Value = 919.3333
$ws.Range("K35").Value = 919.3333
$ws.Range("M35").Value = -540.3333

$ws.Range("H64").Value = 8808.652
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 9266.619000000001
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 9266.619000000001
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -9762.619000000001

$ws.Range("H67").Value = 8808.652
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 9266.619000000001
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 9266.619000000001
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -10982.619

$ws.Range("H138").Value = 1192.8108
$ws.Range("I138").Value = 1132.1945
$ws.Range("K138").Value = 3396.5835
$ws.Range("M138").Value = 1743.4165

$ws.Range("H141").Value = 3851.9375
$ws.Range("I141").Value = 3325.7778
$ws.Range("J141").Value = 4528.4287
$ws.Range("K141").Value = 9977.3334
$ws.Range("L141").Value = 13585.2861
$ws.Range("M141").Value = -4797.3334
$ws.Range("N141").Value = -23945.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1881.6849
$ws.Range("I32").Value = 1900.875
$ws.Range("K32").Value = 1900.875
$ws.Range("M32").Value = -1613.875

$ws.Range("H36").Value = 9206.75
$ws.Range("I36").Value = 9206.75
$ws.Range("K36").Value = 9206.75
$ws.Range("M36").Value = -8860.75

$ws.Range("H63").Value = 3826.6667
$ws.Range("I63").Value = 1391
$ws.Range("J63").Value = 7236.6
$ws.Range("K63").Value = 1391
$ws.Range("L63").Value = 7236.6
$ws.Range("M63").Value = -705
$ws.Range("N63").Value = -8608.6

$ws.Range("H66").Value = 3826.6667
$ws.Range("I66").Value = 1391
$ws.Range("J66").Value = 7236.6
$ws.Range("K66").Value = 6955
$ws.Range("L66").Value = 36183
$ws.Range("M66").Value = -3523
$ws.Range("N66").Value = -43047

$ws.Range("H103").Value = 20000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 20000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 20000
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -22344

$ws.Range("H122").Value = 8548377
$ws.Range("I122").Value = 9260491
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 27781473
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -27779023
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 2126.162
$ws.Range("I132").Value = 2113.543
$ws.Range("K132").Value = 6340.629000000001
$ws.Range("M132").Value = -3810.629000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 16500
$ws.Range("I33").Value = 16500
$ws.Range("K33").Value = 16500
$ws.Range("M33").Value = -16164

$ws.Range("H134").Value = 2319.3572
$ws.Range("I134").Value = 2214.25
$ws.Range("K134").Value = 6642.75
$ws.Range("M134").Value = -4107.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 328.06668
$ws.Range("I7").Value = 208.66667
$ws.Range("J7").Value = 507.16666
$ws.Range("K7").Value = 208.66667
$ws.Range("L7").Value = 507.16666
$ws.Range("M7").Value = -95.66667000000001
$ws.Range("N7").Value = -733.16666

$ws.Range("H22").Value = 1364.6666
$ws.Range("J22").Value = 1258
$ws.Range("L22").Value = 1258
$ws.Range("N22").Value = -1958

$ws.Range("H31").Value = 3227.2354
$ws.Range("I31").Value = 1258.9231
$ws.Range("J31").Value = 9624.25
$ws.Range("K31").Value = 1258.9231
$ws.Range("L31").Value = 9624.25
$ws.Range("M31").Value = -963.9231
$ws.Range("N31").Value = -10214.25

$ws.Range("H34").Value = 3227.2354
$ws.Range("I34").Value = 1258.9231
$ws.Range("J34").Value = 9624.25
$ws.Range("K34").Value = 1258.9231
$ws.Range("L34").Value = 9624.25
$ws.Range("M34").Value = -1056.9231
$ws.Range("N34").Value = -10028.25

$ws.Range("H99").Value = 1115047.4
$ws.Range("I99").Value = 1433104.8
$ws.Range("J99").Value = 1846.5
$ws.Range("K99").Value = 1433104.8
$ws.Range("L99").Value = 1846.5
$ws.Range("M99").Value = -1431606.8
$ws.Range("N99").Value = -4842.5

$ws.Range("H119").Value = 72210
$ws.Range("J119").Value = 72210
$ws.Range("L119").Value = 72210
$ws.Range("N119").Value = -81886

$ws.Range("H122").Value = 1032.4546
$ws.Range("I122").Value = 915
$ws.Range("K122").Value = 2745
$ws.Range("M122").Value = -295

$ws.Range("H126").Value = 1115047.4
$ws.Range("I126").Value = 1433104.8
$ws.Range("J126").Value = 1846.5
$ws.Range("K126").Value = 4299314.4
$ws.Range("L126").Value = 5539.5
$ws.Range("M126").Value = -4296844.4
$ws.Range("N126").Value = -10479.5

$ws.Range("H140").Value = 270170
$ws.Range("J140").Value = 276893.34
$ws.Range("L140").Value = 276893.34
$ws.Range("N140").Value = -287253.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2499.75
$ws.Range("I132").Value = 2049.6667
$ws.Range("J132").Value = 2769.8
$ws.Range("K132").Value = 18447.0003
$ws.Range("L132").Value = 24928.2
$ws.Range("M132").Value = -15917.0003
$ws.Range("N132").Value = -29988.2

$ws.Range("H140").Value = 1629.9166
$ws.Range("I140").Value = 1629.9166
$ws.Range("K140").Value = 4889.7498
$ws.Range("M140").Value = 290.2502000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 26249.75
$ws.Range("J49").Value = 26249.75
$ws.Range("L49").Value = 26249.75
$ws.Range("N49").Value = -26617.75

$ws.Range("H107").Value = 1147.3636
$ws.Range("J107").Value = 1575
$ws.Range("L107").Value = 1575
$ws.Range("N107").Value = -5415

$ws.Range("H122").Value = 17555.842
$ws.Range("I122").Value = 18560.268
$ws.Range("J122").Value = 13789.25
$ws.Range("K122").Value = 55680.804
$ws.Range("L122").Value = 41367.75
$ws.Range("M122").Value = -53230.804
$ws.Range("N122").Value = -46267.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3036.5925
$ws.Range("I22").Value = 2953
$ws.Range("J22").Value = 3085.7646
$ws.Range("K22").Value = 2953
$ws.Range("L22").Value = 3085.7646
$ws.Range("M22").Value = -2658
$ws.Range("N22").Value = -3675.7646

$ws.Range("H27").Value = 3036.5925
$ws.Range("I27").Value = 2953
$ws.Range("J27").Value = 3085.7646
$ws.Range("K27").Value = 2953
$ws.Range("L27").Value = 3085.7646
$ws.Range("M27").Value = -2846
$ws.Range("N27").Value = -3299.7646

$ws.Range("H33").Value = 18960
$ws.Range("J33").Value = 18960
$ws.Range("L33").Value = 18960
$ws.Range("N33").Value = -19540

$ws.Range("H43").Value = 22325.37
$ws.Range("I43").Value = 21499.066
$ws.Range("J43").Value = 23874.688
$ws.Range("K43").Value = 21499.066
$ws.Range("L43").Value = 23874.688
$ws.Range("M43").Value = -21306.066
$ws.Range("N43").Value = -24260.688

$ws.Range("H55").Value = 2005.0769
$ws.Range("I55").Value = 262.16666
$ws.Range("J55").Value = 3499
$ws.Range("K55").Value = 262.16666
$ws.Range("L55").Value = 3499
$ws.Range("M55").Value = -89.16665999999998
$ws.Range("N55").Value = -3845

$ws.Range("H122").Value = 6036.615
$ws.Range("I122").Value = 5497.8184
$ws.Range("K122").Value = 16493.4552
$ws.Range("M122").Value = -14043.4552

$ws.Range("H136").Value = 11848.73
$ws.Range("I136").Value = 1314.1428
$ws.Range("J136").Value = 15729.895
$ws.Range("K136").Value = 3942.4284
$ws.Range("L136").Value = 47189.685
$ws.Range("M136").Value = -1392.4284
$ws.Range("N136").Value = -52289.685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 49750
$ws.Range("J98").Value = 49750
$ws.Range("L98").Value = 49750
$ws.Range("N98").Value = -55740

$ws.Range("H122").Value = 5200.5
$ws.Range("I122").Value = 2355.889
$ws.Range("J122").Value = 8857.857
$ws.Range("K122").Value = 7067.667
$ws.Range("L122").Value = 26573.571
$ws.Range("M122").Value = -4617.667
$ws.Range("N122").Value = -31473.571

$ws.Range("H126").Value = 2353.7
$ws.Range("I126").Value = 1857.6
$ws.Range("J126").Value = 2849.8
$ws.Range("K126").Value = 5572.799999999999
$ws.Range("L126").Value = 8549.400000000001
$ws.Range("M126").Value = -3102.799999999999
$ws.Range("N126").Value = -13489.4

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
